# Applies the "Added multiThreading and Multiprocessing" edit:
#  - Clears the stray fold-id values (col F) that were left on the
#    non-aggregate rows of each 5-row block in Results!A37:I100.
#  - Removes the now-unused scratch/summary rows 104:120 on the Results
#    sheet (raw pivot-style numbers that lived below the real table).
#  - Updates the active sheet / selection so "Time" is the active tab
#    (matches the workbook's stored view state after the edit).

$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("Results")
$wsTime = $wb.Worksheets.Item("Time")

# --- Clear column F on the non-aggregate rows of each block ---------------
$blockStarts = 37, 47, 57, 67, 77, 87, 97
foreach ($start in $blockStarts) {
    for ($r = $start; $r -le ($start + 3); $r++) {
        $wsResults.Range("F$r").ClearContents()
    }
}

# --- Drop the leftover scratch rows below the real table -------------------
$wsResults.Range("A104:U120").EntireRow.Delete()

# --- Restore the stored selections on each sheet ---------------------------
[void]$wsResults.Range("I6:I96").Select()
[void]$wsTime.Range("A2:A33").Select()

# --- Make "Time" the active sheet/tab --------------------------------------
[void]$wsTime.Activate()
